# Updates Price (D) and Volume(1h) (E) columns in the cryptos list
# to reflect refreshed data, per the GitHub Actions scheduled update.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

function Set-TextValue {
    param($cellRef, $val)
    # Force the value to be stored as text (matching the original
    # inline-string cells) instead of letting Excel auto-convert
    # numeric-looking strings into numbers, then restore the default
    # (unstyled) cell style so no extra formatting/quote-prefix marker
    # is introduced.
    $rng = $ws.Range($cellRef)
    $rng.NumberFormat = "@"
    $rng.Value = $val
    $rng.Style = "Normal"
}

Set-TextValue 'D2' '64.379.03'
Set-TextValue 'E2' '  -3.64%  '
Set-TextValue 'D3' '3.158.50'
Set-TextValue 'E3' '  -3.15%  '
Set-TextValue 'E4' '  +0.22%  '
Set-TextValue 'D5' '608.19'
Set-TextValue 'E5' '  -0.14%  '
Set-TextValue 'D6' '147.72'
Set-TextValue 'E6' '  -6.84%  '
Set-TextValue 'E7' '  +0.14%  '
Set-TextValue 'D8' '3.151.96'
Set-TextValue 'E8' '  -3.41%  '
Set-TextValue 'D9' '0.525'
Set-TextValue 'E9' '  -4.20%  '
Set-TextValue 'E10' '  -6.56%  '
Set-TextValue 'D11' '5.49'
Set-TextValue 'E11' '  -6.94%  '
Set-TextValue 'E12' '  -5.97%  '
Set-TextValue 'E13' '  -7.94%  '
Set-TextValue 'D14' '35.64'
Set-TextValue 'E14' '  -9.42%  '
Set-TextValue 'D15' '3.677.60'
Set-TextValue 'E15' '  -2.98%  '
Set-TextValue 'D16' '64.390.29'
Set-TextValue 'E16' '  -3.61%  '
Set-TextValue 'D18' '3.157.73'
Set-TextValue 'E18' '  -3.97%  '
Set-TextValue 'E19' '  -6.81%  '
Set-TextValue 'D20' '481.47'
Set-TextValue 'E20' '  -5.57%  '
Set-TextValue 'D21' '14.74'
Set-TextValue 'E21' '  -4.68%  '
Set-TextValue 'E22' '  -5.69%  '
Set-TextValue 'D23' '7.79'
Set-TextValue 'D24' '13.71'
Set-TextValue 'E24' '  -7.98%  '
Set-TextValue 'D25' '83.73'
Set-TextValue 'E25' '  -3.39%  '
Set-TextValue 'E26' '  -0.07%  '
Set-TextValue 'E27' '  -5.36%  '
Set-TextValue 'E28' '  -7.59%  '
Set-TextValue 'E29' '  -9.42%  '
Set-TextValue 'D30' '6.83'
Set-TextValue 'E30' '  -0.76%  '
Set-TextValue 'E31' '  -19.57%  '
Set-TextValue 'E32' '  -5.90%  '
Set-TextValue 'E33' '  +0.10%  '
Set-TextValue 'D34' '26.26'
Set-TextValue 'E34' '  -6.82%  '
Set-TextValue 'E35' '  -4.72%  '
Set-TextValue 'D36' '54.54'
Set-TextValue 'E36' '  -2.41%  '
Set-TextValue 'D37' '5.98'
Set-TextValue 'E37' '  -7.69%  '
Set-TextValue 'E38' '  -9.01%  '
Set-TextValue 'D39' '458.14'
Set-TextValue 'E39' '  -7.92%  '
Set-TextValue 'E40' '  -14.05%  '
Set-TextValue 'E41' '  -7.46%  '
Set-TextValue 'E42' '  -4.86%  '
Set-TextValue 'E43' '  -8.23%  '
Set-TextValue 'D44' '2.855.67'
Set-TextValue 'E44' '  -4.11%  '
Set-TextValue 'E45' '  -10.02%  '
Set-TextValue 'E46' '  -9.62%  '
Set-TextValue 'D47' '26.59'
Set-TextValue 'E47' '  -8.28%  '
Set-TextValue 'E48' '  +0.02%  '
Set-TextValue 'E49' '  -7.42%  '
Set-TextValue 'E50' '  -4.85%  '
Set-TextValue 'D51' '119.85'
Set-TextValue 'E51' '  -1.55%  '
